$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was inserted as row 260 ("Hortaliza, Terminal
# Hortofrutícola Agro Chillán - Pepino ensalada"), pushing the existing rows
# 260-309 down to 261-310 (same data, unchanged).
$ws.Rows.Item(260).Insert()

$ws.Cells.Item(260, 1).Value = 7
$ws.Cells.Item(260, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(260, 3).Value = "Ñuble"
$ws.Cells.Item(260, 4).Value = 45015
$ws.Cells.Item(260, 5).Value = 16
$ws.Cells.Item(260, 6).Value = 100112043
$ws.Cells.Item(260, 7).Value = "Pepino ensalada"
$ws.Cells.Item(260, 8).Value = "Sin especificar"
$ws.Cells.Item(260, 9).Value = "Primera"
$ws.Cells.Item(260, 10).Value = 60
$ws.Cells.Item(260, 11).Value = 13000
$ws.Cells.Item(260, 12).Value = 13000
$ws.Cells.Item(260, 13).Value = 13000
$ws.Cells.Item(260, 14).Value = "`$/caja 80 unidades"
$ws.Cells.Item(260, 15).Value = "Región del Maule"
$ws.Cells.Item(260, 16).Value = 162
$ws.Cells.Item(260, 17).Value = 80
$ws.Cells.Item(260, 18).Value = "Hortaliza"
